# workplan.xlsx update — see commit message "updated the workplan file"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix / rename existing cell text -----------------------------------
$ws.Range("A2").Value  = "Youth Demographics"
$ws.Range("A4").Value  = "Youth Demographics"
$ws.Range("A5").Value  = "Youth Demographics"
$ws.Range("A6").Value  = "pyschometrics "

# --- new "comment" column (C) ------------------------------------------
$ws.Range("C1").Value = "comment"
$ws.Range("C4").Value = "done"

# --- new rows appended below the existing table -------------------------
$ws.Range("A7").Value  = "pyschometrics "
$ws.Range("B7").Value  = "add weekly and monthly windows"

$ws.Range("A8").Value  = "Youth Demographics"
$ws.Range("B8").Value  = "add placed summary on top left"

$ws.Range("A9").Value  = "Engagements"
$ws.Range("B9").Value  = "add top performing companies"
$ws.Range("C9").Value  = " {how do we say a company is top / what do we consider}"

$ws.Range("A10").Value = "Engagements"
$ws.Range("B10").Value = "add top performing youth {high answered/expected}"
$ws.Range("C10").Value = "do we have a minimum placed time"

$ws.Range("A11").Value = "pyschometrics "
$ws.Range("B11").Value = "Weekly and Monthly histograms"

# --- header row formatting: double bottom border, then bold ------------
$headerRange = $ws.Range("A1:C1")
$headerRange.Borders.Item(9).LineStyle = -4119   # xlDouble
$headerRange.Font.Bold = $true

# --- column widths (auto best-fit to new content) -----------------------
$ws.Columns("A:C").AutoFit()

# --- page setup: portrait orientation -----------------------------------
$ws.PageSetup.Orientation = 1   # xlPortrait

# --- leave the cursor where the author left it when saving --------------
$ws.Range("B14").Select() | Out-Null
